$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.847.03'
$ws.Range('E2').Value = '  +2.31%  '
$ws.Range('D3').Value = '2.337.07'
$ws.Range('E3').Value = '  +2.29%  '
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range('D5').Value = '312.01'
$ws.Range('E5').Value = '  -0.50%  '
$ws.Range('D6').Value = '108.54'
$ws.Range('E6').Value = '  +3.87%  '
$ws.Range('D7').Value = '0.633'
$ws.Range('E7').Value = '  +1.09%  '
$ws.Range('E8').Value = '  -0.16%  '
$ws.Range('D9').Value = '0.619'
$ws.Range('E9').Value = '  +2.34%  '
$ws.Range('D10').Value = '41.22'
$ws.Range('E10').Value = '  +4.66%  '
$ws.Range('D11').Value = '0.0919'
$ws.Range('E11').Value = '  +1.68%  '
$ws.Range('D12').Value = '8.55'
$ws.Range('E12').Value = '  +2.34%  '
$ws.Range('E13').Value = '  -1.36%  '
$ws.Range('D14').Value = '1.01'
$ws.Range('E14').Value = '  +1.62%  '
$ws.Range('D15').Value = '15.46'
$ws.Range('E15').Value = '  +1.54%  '
$ws.Range('D16').Value = '2.691.32'
$ws.Range('E16').Value = '  +2.13%  '
$ws.Range('D17').Value = '2.334.40'
$ws.Range('E17').Value = '  +2.24%  '
$ws.Range('D18').Value = '43.759.61'
$ws.Range('E18').Value = '  +2.31%  '
$ws.Range('D19').Value = '7.55'
$ws.Range('E19').Value = '  +1.91%  '
$ws.Range('E20').Value = '  +1.33%  '
$ws.Range('E21').Value = '  -3.06%  '
$ws.Range('D22').Value = '74.21'
$ws.Range('E22').Value = '  +0.45%  '
$ws.Range('D23').Value = '3.47'
$ws.Range('E23').Value = '  -3.38%  '
$ws.Range('D24').Value = '268.88'
$ws.Range('E24').Value = '  +1.49%  '
$ws.Range('D25').Value = '2.28'
$ws.Range('E25').Value = '  +2.91%  '
$ws.Range('E26').Value = '  -0.11%  '
$ws.Range('E27').Value = '  +7.23%  '
$ws.Range('D28').Value = '11.13'
$ws.Range('E28').Value = '  +2.83%  '
$ws.Range('E29').Value = '  -1.90%  '
$ws.Range('D30').Value = '39.25'
$ws.Range('E30').Value = '  +5.74%  '
$ws.Range('D31').Value = '22.61'
$ws.Range('E31').Value = '  +0.60%  '
$ws.Range('D32').Value = '168.57'
$ws.Range('E32').Value = '  +0.88%  '
$ws.Range('D33').Value = '0.0887'
$ws.Range('E33').Value = '  +1.65%  '
$ws.Range('D34').Value = '2.84'
$ws.Range('E34').Value = '  +9.52%  '
$ws.Range('E35').Value = '  +0.85%  '
$ws.Range('D36').Value = '0.115'
$ws.Range('E36').Value = '  +0.92%  '
$ws.Range('E37').Value = '  +3.73%  '
$ws.Range('E38').Value = '  +3.26%  '
$ws.Range('D39').Value = '2.89'
$ws.Range('E39').Value = '  +8.45%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '3.80'
$ws.Range('E40').Value = '  +1.01%  '
$ws.Range('D41').Value = '1.71'
$ws.Range('E41').Value = '  +8.47%  '
$ws.Range('D42').Value = '104.87'
$ws.Range('E42').Value = '  +11.75%  '
$ws.Range('D43').Value = '0.239'
$ws.Range('E43').Value = '  +2.79%  '
$ws.Range('D44').Value = '13.47'
$ws.Range('E44').Value = '  +11.27%  '
$ws.Range('D45').Value = '71.76'
$ws.Range('E45').Value = '  +1.65%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '1.00'
$ws.Range('E46').Value = '  -0.07%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '114.00'
$ws.Range('E47').Value = '  +1.18%  '
$ws.Range('D48').Value = '1.669.14'
$ws.Range('E48').Value = '  -4.20%  '
$ws.Range('D49').Value = '0.219'
$ws.Range('E49').Value = '  +16.11%  '
$ws.Range('D50').Value = '76.92'
$ws.Range('E50').Value = '  -3.54%  '
$ws.Range('D51').Value = '8.95'
$ws.Range('E51').Value = '  +2.48%  '
